# Sales Invoice workbook update: add Tax / Customer Notes / Terms And Conditions /
# Save As / Price List columns to the header sheet, add Discount Type / Discount
# columns to the items sheet, and refresh the sample date values.

$wb = $excel.ActiveWorkbook

$wsHeader = $wb.Worksheets.Item("SalesInvoiceHeader")
$wsItems  = $wb.Worksheets.Item("SalesInvoiceItems")

# --- SalesInvoiceHeader ---------------------------------------------------

# New header labels (row 1, columns I:M)
$wsHeader.Cells.Item(1, 9).Value  = "Tax"
$wsHeader.Cells.Item(1, 10).Value = "Customer Notes"
$wsHeader.Cells.Item(1, 11).Value = "Terms And Conditions"
$wsHeader.Cells.Item(1, 12).Value = "Save As"
$wsHeader.Cells.Item(1, 13).Value = "Price List"

# New sample values (row 2, columns I:M)
$wsHeader.Cells.Item(2, 9).Value  = "Inclusive"
$wsHeader.Cells.Item(2, 10).Value = "notex"
$wsHeader.Cells.Item(2, 11).Value = "termsx"
$wsHeader.Cells.Item(2, 12).Value = "SAVE AND APPROVE"
$wsHeader.Cells.Item(2, 13).Value = "special price"

# Updated sample dates
$wsHeader.Cells.Item(2, 2).Value = "271225"
$wsHeader.Cells.Item(2, 3).Value = "27-12-2025"
$wsHeader.Cells.Item(2, 5).Value = "30-12-2025"

# --- SalesInvoiceItems -----------------------------------------------------

$wsItems.Cells.Item(1, 4).Value = "Discount Type"
$wsItems.Cells.Item(1, 5).Value = "Discount"

$wsItems.Cells.Item(2, 4).Value = "%"
$wsItems.Cells.Item(2, 5).Value = 10

$wsItems.Cells.Item(3, 4).Value = "amount"
$wsItems.Cells.Item(3, 5).Value = 5

# --- Active sheet / selection ----------------------------------------------
# The workbook now opens on the Items sheet rather than the Header sheet.
$wsItems.Select()
$wsHeader.Range("F6").Select()
$wsItems.Range("Q5").Select()
